$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sending-cluster set now includes "ECs" in addition to "FAPs" and "MuSCs",
# so the data grid grows from 2 sending clusters x 3 target clusters (rows 2-7)
# to 3 sending clusters x 3 target clusters (rows 2-10), and all the dependent
# specificity/weight figures are recalculated against the larger group.

$data = @{
  2 = @("ECs", "Efnb3", "Ephb6", "ECs", 2, 0.6666666666666666, 0.2151756666666667, 0.645527, 0.1791915537270423, 0.1791915537270423, 2, 0.6666666666666666, 0.08324533333333334, 0.249736, 0.05078606388889115, 0.05078606388889115, 0.01791237009688889, 0.161211330872, 0.009100433695931241, 0.009100433695931241)
  3 = @("ECs", "Efnb3", "Ephb6", "FAPs", 2, 0.6666666666666666, 0.2151756666666667, 0.645527, 0.1791915537270423, 0.1791915537270423, 2, 0.6666666666666666, 0.4583163333333333, 1.374949, 0.2796082573516313, 0.2796082573516313, 0.09861852256922221, 0.8875667031229999, 0.05010343806974951, 0.0501034380697495)
  4 = @("ECs", "Efnb3", "Ephb6", "MuSCs", 2, 0.6666666666666666, 0.2151756666666667, 0.645527, 0.1791915537270423, 0.1791915537270423, 3, 1, 1.097575666666667, 3.292727, 0.6696056787594775, 0.6696056787594775, 0.2361715757921111, 2.125544182129, 0.1199876819613615, 0.1199876819613615)
  5 = @("FAPs", "Efnb3", "Ephb6", "ECs", 2, 0.6666666666666666, 0.5010686666666667, 1.503206, 0.4172742870736815, 0.4172742870736815, 2, 0.6666666666666666, 0.08324533333333334, 0.249736, 0.05078606388889115, 0.05078606388889115, 0.04171162817955557, 0.375404653616, 0.02119171860251549, 0.02119171860251549)
  6 = @("FAPs", "Efnb3", "Ephb6", "FAPs", 2, 0.6666666666666666, 0.5010686666666667, 1.503206, 0.4172742870736815, 0.4172742870736815, 2, 0.6666666666666666, 0.4583163333333333, 1.374949, 0.2796082573516313, 0.2796082573516313, 0.2296479540548889, 2.066831586494, 0.1166733362463164, 0.1166733362463164)
  7 = @("FAPs", "Efnb3", "Ephb6", "MuSCs", 2, 0.6666666666666666, 0.5010686666666667, 1.503206, 0.4172742870736815, 0.4172742870736815, 3, 1, 1.097575666666667, 3.292727, 0.6696056787594775, 0.6696056787594775, 0.5499607758624445, 4.949646982762, 0.2794092322248495, 0.2794092322248495)
  8 = @("MuSCs", "Efnb3", "Ephb6", "ECs", 2, 0.6666666666666666, 0.4845693333333334, 1.453708, 0.4035341591992763, 0.4035341591992763, 2, 0.6666666666666666, 0.08324533333333334, 0.249736, 0.05078606388889115, 0.05078606388889115, 0.04033813567644445, 0.363043221088, 0.02049391159044442, 0.02049391159044442)
  9 = @("MuSCs", "Efnb3", "Ephb6", "FAPs", 2, 0.6666666666666666, 0.4845693333333334, 1.453708, 0.4035341591992763, 0.4035341591992763, 2, 0.6666666666666666, 0.4583163333333333, 1.374949, 0.2796082573516313, 0.2796082573516313, 0.2220860400991111, 1.998774360892, 0.1128314830355654, 0.1128314830355654)
  10 = @("MuSCs", "Efnb3", "Ephb6", "MuSCs", 2, 0.6666666666666666, 0.4845693333333334, 1.453708, 0.4035341591992763, 0.4035341591992763, 3, 1, 1.097575666666667, 3.292727, 0.6696056787594775, 0.6696056787594775, 0.5318515090795556, 4.786663581716001, 0.2702087645732664, 0.2702087645732664)
}

foreach ($r in $data.Keys) {
  $rowVals = $data[$r]
  for ($i = 0; $i -lt $rowVals.Count; $i++) {
    $ws.Cells.Item([int]$r, $i + 1).Value = $rowVals[$i]
  }
}
